# feat: add 2022-Q1 data
#
# Before:  2021-Q3, 2021-Q4, 总计
# After:   2021-Q3, 2021-Q4, 2022-Q1, 总计
#
# The old "总计" worksheet (sheetId 3) is repurposed/renamed into the new
# "2022-Q1" fund-holdings sheet (it keeps sheetId 3, matching how Excel
# reuses the id of a renamed sheet), and a brand-new "总计" worksheet is
# appended after it (getting a fresh sheetId 4) with the refreshed totals
# table that now also lists the 2022-Q1 quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Turn the old "总计" sheet into the new "2022-Q1" holdings sheet.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# Pull the header formatting (bold font + border + centered) from the
# 2021-Q4 sheet's header row so the new columns (E:H) match the look of
# every other quarter sheet, without inventing a brand-new style.
$prevHeader = $wb.Worksheets.Item("2021-Q4").Range("B1:H1")
$prevHeader.Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Row 2 - 012442 / 永赢稳健增长一年持有期混合E
$q1.Range("B2").NumberFormat = "@"
$q1.Range("B2").Value = "012442"
$q1.Range("B2").Style = "Normal"
$q1.Range("C2").Value = "永赢稳健增长一年持有期混合E"
$q1.Range("D2").NumberFormat = "@"
$q1.Range("D2").Value = "22.52"
$q1.Range("D2").Style = "Normal"
$q1.Range("E2").NumberFormat = "@"
$q1.Range("E2").Value = "22.68"
$q1.Range("E2").Style = "Normal"
$q1.Range("F2").NumberFormat = "@"
$q1.Range("F2").Value = "2.31"
$q1.Range("F2").Style = "Normal"
$q1.Range("G2").NumberFormat = "@"
$q1.Range("G2").Value = "0.5202"
$q1.Range("G2").Style = "Normal"
$q1.Range("H2").Value = 2

# Row 3 - 009932 / 永赢稳健增长一年持有期混合
$q1.Range("B3").NumberFormat = "@"
$q1.Range("B3").Value = "009932"
$q1.Range("B3").Style = "Normal"
$q1.Range("C3").Value = "永赢稳健增长一年持有期混合"
$q1.Range("D3").NumberFormat = "@"
$q1.Range("D3").Value = "22.52"
$q1.Range("D3").Style = "Normal"
$q1.Range("E3").NumberFormat = "@"
$q1.Range("E3").Value = "22.68"
$q1.Range("E3").Style = "Normal"
$q1.Range("F3").NumberFormat = "@"
$q1.Range("F3").Value = "2.31"
$q1.Range("F3").Style = "Normal"
$q1.Range("G3").NumberFormat = "@"
$q1.Range("G3").Value = "0.5202"
$q1.Range("G3").Style = "Normal"
$q1.Range("H3").Value = 2

# ---------------------------------------------------------------------
# 2) Append the fresh "总计" sheet (right after "2022-Q1") with the
#    updated totals table (2022-Q1 row inserted on top).
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$q1.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$q1.Range("A2").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 1.04

$q1.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 4
$total.Range("D3").Value = 0.96

$q1.Range("A2").Copy()
$total.Range("A4").PasteSpecial(-4122)
$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 5
$total.Range("D4").Value = 0.98
